$d = $word.ActiveDocument

# --- Rename the "_Hlk117032442" bookmark to "_GoBack", and drop the stray
#     one-letter "к" run it wraps (right before "Міністерство освіти і
#     науки України " in the document's very first paragraph).
#
# The bookmark currently spans exactly that one-character run. We capture
# its range, add the new "_GoBack" bookmark around the same (non-empty)
# run, delete the old bookmark, and only then delete the "к" text itself -
# the now-empty "_GoBack" bookmark collapses in place, right before
# "Міністерство". (Adding a bookmark straight onto an empty range sitting
# at document position 0 hits an edge-case in this host, so we go through
# the one-character range first and delete the text afterwards instead.)
$oldBm = $d.Bookmarks("_Hlk117032442")
$start = $oldBm.Start
$letterRange = $d.Range($start, $start + 1)

$oldBm.Delete()
$d.Bookmarks.Add("_GoBack", $letterRange)

$d.Range($start, $start + 1).Delete()

# Word only allows one bookmark per name, so the document's other,
# later "_GoBack" bookmark (around the empty space right after
# "Лектор" + "ка", before ":") is automatically removed by adding the
# "_GoBack" bookmark above - matching the target edit, which drops that
# duplicate bookmarkStart/bookmarkEnd pair entirely.
